$wb = $excel.ActiveWorkbook

# --- Text change: "Ready for handoff" -> "In Translation" ---------------
# Appears on the Overview sheet (E2, F2) and on each per-locale status
# sheet (column C, row 2).
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = "In Translation"
$wsOverview.Range("F2").Value = "In Translation"

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C2").Value = "In Translation"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C2").Value = "In Translation"

# --- Column width change: narrower "Status" columns ----------------------
# The stored width shrinks from ~17.216 to ~13.410 characters. The COM
# ColumnWidth setter here snaps internally to the nearest 1/6th of a
# character, so 12.5 is the input that lands in the middle of the bucket
# that serializes to the target width.
$wsOverview.Range("E1:F1").ColumnWidth = 12.5
$wsZhCn.Range("C1").ColumnWidth = 12.5
$wsDeDe.Range("C1").ColumnWidth = 12.5
